$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 924.75
$ws.Range("I46").Value = 750
$ws.Range("J46").Value = 1099.5
$ws.Range("K46").Value = 2250
$ws.Range("L46").Value = 3298.5
$ws.Range("M46").Value = -2131
$ws.Range("N46").Value = -3536.5

$ws.Range("H60").Value = 924.75
$ws.Range("I60").Value = 750
$ws.Range("J60").Value = 1099.5
$ws.Range("K60").Value = 2250
$ws.Range("L60").Value = 3298.5
$ws.Range("M60").Value = -1766
$ws.Range("N60").Value = -4266.5

$ws.Range("H64").Value = 3257
$ws.Range("I64").Value = 3159.8
$ws.Range("J64").Value = 3500
$ws.Range("K64").Value = 3159.8
$ws.Range("L64").Value = 3500
$ws.Range("M64").Value = -2911.8
$ws.Range("N64").Value = -3996

$ws.Range("H67").Value = 3257
$ws.Range("I67").Value = 3159.8
$ws.Range("J67").Value = 3500
$ws.Range("K67").Value = 3159.8
$ws.Range("L67").Value = 3500
$ws.Range("M67").Value = -2301.8
$ws.Range("N67").Value = -5216

$ws.Range("H76").Value = 2926640.2
$ws.Range("I76").Value = 2708.2856
$ws.Range("J76").Value = 4632267
$ws.Range("K76").Value = 2708.2856
$ws.Range("L76").Value = 4632267
$ws.Range("M76").Value = -2393.2856
$ws.Range("N76").Value = -4632897

$ws.Range("H79").Value = 2926640.2
$ws.Range("I79").Value = 2708.2856
$ws.Range("J79").Value = 4632267
$ws.Range("K79").Value = 2708.2856
$ws.Range("L79").Value = 4632267
$ws.Range("M79").Value = -1616.2856
$ws.Range("N79").Value = -4634451

$ws.Range("H116").Value = 16670950
$ws.Range("I116").Value = 83334830
$ws.Range("J116").Value = 4979.6665
$ws.Range("K116").Value = 83334830
$ws.Range("L116").Value = 4979.6665
$ws.Range("M116").Value = -83331388
$ws.Range("N116").Value = -11863.6665

$ws.Range("H129").Value = 1901.439
$ws.Range("I129").Value = 199
$ws.Range("J129").Value = 2035.8422
$ws.Range("K129").Value = 597
$ws.Range("L129").Value = 6107.5266
$ws.Range("M129").Value = 4403
$ws.Range("N129").Value = -16107.5266

$ws.Range("H141").Value = 1528.9762
$ws.Range("I141").Value = 1310.8649
$ws.Range("J141").Value = 3143
$ws.Range("K141").Value = 3932.5947
$ws.Range("L141").Value = 9429
$ws.Range("M141").Value = 1247.4053
$ws.Range("N141").Value = -19789

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1352.6666
$ws.Range("I2").Value = 1335.4546
$ws.Range("J2").Value = 1400
$ws.Range("K2").Value = 1335.4546
$ws.Range("L2").Value = 1400
$ws.Range("M2").Value = -1222.4546
$ws.Range("N2").Value = -1626

$ws.Range("H32").Value = 8522.937
$ws.Range("I32").Value = 6106.5737
$ws.Range("J32").Value = 23460.455
$ws.Range("K32").Value = 6106.5737
$ws.Range("L32").Value = 23460.455
$ws.Range("M32").Value = -5819.5737
$ws.Range("N32").Value = -24034.455

$ws.Range("H63").Value = 2843566.2
$ws.Range("I63").Value = 2914.2222
$ws.Range("J63").Value = 15626500
$ws.Range("K63").Value = 2914.2222
$ws.Range("L63").Value = 15626500
$ws.Range("M63").Value = -2228.2222
$ws.Range("N63").Value = -15627872

$ws.Range("H66").Value = 2843566.2
$ws.Range("I66").Value = 2914.2222
$ws.Range("J66").Value = 15626500
$ws.Range("K66").Value = 14571.111
$ws.Range("L66").Value = 78132500
$ws.Range("M66").Value = -11139.111
$ws.Range("N66").Value = -78139364

$ws.Range("H102").Value = 2010
$ws.Range("I102").Value = 2010
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2010
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -388

$ws.Range("H116").Value = 1352.6666
$ws.Range("I116").Value = 1335.4546
$ws.Range("J116").Value = 1400
$ws.Range("K116").Value = 1335.4546
$ws.Range("L116").Value = 1400
$ws.Range("M116").Value = 958.5454
$ws.Range("N116").Value = -5988

$ws.Range("H122").Value = 1577.9788
$ws.Range("I122").Value = 1501.45
$ws.Range("J122").Value = 2015.2858
$ws.Range("K122").Value = 4504.35
$ws.Range("L122").Value = 6045.857400000001
$ws.Range("M122").Value = -2054.35
$ws.Range("N122").Value = -10945.8574

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1352.6666
$ws.Range("I3").Value = 1335.4546
$ws.Range("J3").Value = 1400
$ws.Range("K3").Value = 1335.4546
$ws.Range("L3").Value = 1400
$ws.Range("M3").Value = -1221.4546
$ws.Range("N3").Value = -1628

$ws.Range("H80").Value = 998.5
$ws.Range("I80").Value = 1263.2727
$ws.Range("J80").Value = 827.17645
$ws.Range("K80").Value = 1263.2727
$ws.Range("L80").Value = 827.17645
$ws.Range("M80").Value = -265.2727
$ws.Range("N80").Value = -2823.17645

$ws.Range("H83").Value = 998.5
$ws.Range("I83").Value = 1263.2727
$ws.Range("J83").Value = 827.17645
$ws.Range("K83").Value = 6316.363499999999
$ws.Range("L83").Value = 4135.882250000001
$ws.Range("M83").Value = -1324.363499999999
$ws.Range("N83").Value = -14119.88225

$ws.Range("H86").Value = 2252.12
$ws.Range("I86").Value = 2262.9473
$ws.Range("J86").Value = 2217.8333
$ws.Range("K86").Value = 2262.9473
$ws.Range("L86").Value = 2217.8333
$ws.Range("M86").Value = -1139.9473
$ws.Range("N86").Value = -4463.8333

$ws.Range("H89").Value = 2252.12
$ws.Range("I89").Value = 2262.9473
$ws.Range("J89").Value = 2217.8333
$ws.Range("K89").Value = 11314.7365
$ws.Range("L89").Value = 11089.1665
$ws.Range("M89").Value = -5698.736499999999
$ws.Range("N89").Value = -22321.1665

$ws.Range("H94").Value = 1407.1428
$ws.Range("I94").Value = 1250
$ws.Range("J94").Value = 1470
$ws.Range("K94").Value = 1250
$ws.Range("L94").Value = 1470
$ws.Range("M94").Value = -799
$ws.Range("N94").Value = -2372

$ws.Range("H105").Value = 1726390.9
$ws.Range("I105").Value = 1608.3334
$ws.Range("J105").Value = 2176334.2
$ws.Range("K105").Value = 1608.3334
$ws.Range("L105").Value = 2176334.2
$ws.Range("M105").Value = 138.6666
$ws.Range("N105").Value = -2179828.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1092.909
$ws.Range("I16").Value = 1080.2222
$ws.Range("J16").Value = 1150
$ws.Range("K16").Value = 1080.2222
$ws.Range("L16").Value = 1150
$ws.Range("M16").Value = -793.2221999999999
$ws.Range("N16").Value = -1724

$ws.Range("H31").Value = 3555.7778
$ws.Range("I31").Value = 1697.9697
$ws.Range("J31").Value = 5599.3667
$ws.Range("K31").Value = 1697.9697
$ws.Range("L31").Value = 5599.3667
$ws.Range("M31").Value = -1402.9697
$ws.Range("N31").Value = -6189.3667

$ws.Range("H34").Value = 3555.7778
$ws.Range("I34").Value = 1697.9697
$ws.Range("J34").Value = 5599.3667
$ws.Range("K34").Value = 1697.9697
$ws.Range("L34").Value = 5599.3667
$ws.Range("M34").Value = -1495.9697
$ws.Range("N34").Value = -6003.3667

$ws.Range("H113").Value = 1092.909
$ws.Range("I113").Value = 1080.2222
$ws.Range("J113").Value = 1150
$ws.Range("K113").Value = 1080.2222
$ws.Range("L113").Value = 1150
$ws.Range("M113").Value = 1089.7778
$ws.Range("N113").Value = -5490

$ws.Range("H132").Value = 2460.3225
$ws.Range("I132").Value = 1688.9166
$ws.Range("J132").Value = 5105.143
$ws.Range("K132").Value = 5066.7498
$ws.Range("L132").Value = 15315.429
$ws.Range("M132").Value = -2536.7498
$ws.Range("N132").Value = -20375.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1239.2122
$ws.Range("I5").Value = 889.2083
$ws.Range("J5").Value = 2172.5557
$ws.Range("K5").Value = 2667.6249
$ws.Range("L5").Value = 6517.6671
$ws.Range("M5").Value = -2555.6249
$ws.Range("N5").Value = -6741.6671

$ws.Range("H8").Value = 345.2
$ws.Range("I8").Value = 345.2
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1035.6
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -896.5999999999999

$ws.Range("H122").Value = 1237.1562
$ws.Range("I122").Value = 375
$ws.Range("J122").Value = 1294.6333
$ws.Range("K122").Value = 3375
$ws.Range("L122").Value = 11651.6997
$ws.Range("M122").Value = -925
$ws.Range("N122").Value = -16551.6997

$ws.Range("H131").Value = 731.5599999999999
$ws.Range("I131").Value = 391.875
$ws.Range("J131").Value = 761.09784
$ws.Range("K131").Value = 1175.625
$ws.Range("L131").Value = 2283.29352
$ws.Range("M131").Value = 3864.375
$ws.Range("N131").Value = -12363.29352

$ws.Range("H135").Value = 1239.2122
$ws.Range("I135").Value = 889.2083
$ws.Range("J135").Value = 2172.5557
$ws.Range("K135").Value = 8002.8747
$ws.Range("L135").Value = 19553.0013
$ws.Range("M135").Value = -5467.8747
$ws.Range("N135").Value = -24623.0013

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3129375.2
$ws.Range("I70").Value = 4599.778
$ws.Range("J70").Value = 5686010
$ws.Range("K70").Value = 4599.778
$ws.Range("L70").Value = 5686010
$ws.Range("M70").Value = -4329.778
$ws.Range("N70").Value = -5686550

$ws.Range("H73").Value = 3129375.2
$ws.Range("I73").Value = 4599.778
$ws.Range("J73").Value = 5686010
$ws.Range("K73").Value = 4599.778
$ws.Range("L73").Value = 5686010
$ws.Range("M73").Value = -3663.778
$ws.Range("N73").Value = -5687882

$ws.Range("H122").Value = 4980.346
$ws.Range("I122").Value = 4852.6665
$ws.Range("J122").Value = 5154.4546
$ws.Range("K122").Value = 14557.9995
$ws.Range("L122").Value = 15463.3638
$ws.Range("M122").Value = -12107.9995
$ws.Range("N122").Value = -20363.3638

$ws.Range("H132").Value = 84777.45
$ws.Range("I132").Value = 104454.45
$ws.Range("J132").Value = 49001.09
$ws.Range("K132").Value = 313363.35
$ws.Range("L132").Value = 147003.27
$ws.Range("M132").Value = -310833.35
$ws.Range("N132").Value = -152063.27

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3942.8572
$ws.Range("I7").Value = 3316.6667
$ws.Range("J7").Value = 7700
$ws.Range("K7").Value = 3316.6667
$ws.Range("L7").Value = 7700
$ws.Range("M7").Value = -3204.6667
$ws.Range("N7").Value = -7924

$ws.Range("H16").Value = 783.6667
$ws.Range("I16").Value = 775.3
$ws.Range("J16").Value = 800.4
$ws.Range("K16").Value = 775.3
$ws.Range("L16").Value = 800.4
$ws.Range("M16").Value = -605.3
$ws.Range("N16").Value = -1140.4

$ws.Range("H126").Value = 3942.8572
$ws.Range("I126").Value = 3316.6667
$ws.Range("J126").Value = 7700
$ws.Range("K126").Value = 9950.000100000001
$ws.Range("L126").Value = 23100
$ws.Range("M126").Value = -7480.000100000001
$ws.Range("N126").Value = -28040

$ws.Range("H132").Value = 213809.4
$ws.Range("I132").Value = 295905.75
$ws.Range("J132").Value = 3437.5
$ws.Range("K132").Value = 887717.25
$ws.Range("L132").Value = 10312.5
$ws.Range("M132").Value = -885187.25
$ws.Range("N132").Value = -15372.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 819.4074000000001
$ws.Range("I113").Value = 1035.6875
$ws.Range("J113").Value = 504.81818
$ws.Range("K113").Value = 3107.0625
$ws.Range("L113").Value = 1514.45454
$ws.Range("M113").Value = -937.0625
$ws.Range("N113").Value = -5854.45454
